# Apply cryptos list update (prices and volume%) per commit:
# "Updated cryptos list on Mon Feb 27 08:33:48 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.414.06"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.636.25"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.27%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.002"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "305.00"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.79%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3740"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.99%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "51.87"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3616"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.252"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08117"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.80"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.599"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001266"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.273"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.636.25"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.23"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.10"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.503"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "23.413.44"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.72"
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.421"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.84%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.047"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.16"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.63"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.317"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "135.78"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.286"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.817.34"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.717"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9512"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02817"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.24"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.07249"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2512"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.08782"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.067"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.371"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7032"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.70%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.11"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.20%  "
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.39"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6502"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.319"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.001"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.007"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07964"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "128.09"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.197"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.06%  "
